$d = $word.ActiveDocument

# --- Fix capitalization / spelling of species epithets ---
$d.Content.Find.Execute("Vivax", $true, $false, $false, $false, $false, $true, 1, $false, "vivax", 2)
$d.Content.Find.Execute("Gondii", $true, $false, $false, $false, $false, $true, 1, $false, "gondii", 2)
$d.Content.Find.Execute("Brucei", $true, $false, $false, $false, $false, $true, 1, $false, "brucei", 2)
$d.Content.Find.Execute("Cruzi", $true, $false, $false, $false, $false, $true, 1, $false, "cruzi", 2)
$d.Content.Find.Execute("Major", $true, $false, $false, $false, $false, $true, 1, $false, "major", 2)
$d.Content.Find.Execute("entamoeba", $true, $false, $false, $false, $false, $true, 1, $false, "Entamoeba", 2)

# --- Move the _GoBack bookmark from after "[BLASTP or HMM" to right after "and E" ---
$r = $d.Content
$r.Find.Execute("and E", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target = $r.End

$d.Bookmarks.Item("_GoBack").Delete()

$newRange = $d.Range($target, $target)
$d.Bookmarks.Add("_GoBack", $newRange)

# --- Re-merge the runs that were split around the old bookmark location ---
$d.Content.Find.Execute("[BLASTP or HMM]", $true, $false, $false, $false, $false, $true, 1, $false, "[BLASTP or HMM]", 2)
